$wb = $excel.ActiveWorkbook

# The e2e/ca6ecd48-c673-4fea-9792-4020c15d8bc0.md handback file finished
# processing: status flips from "Ready for handoff" to
# "Handed back: in sync with en-US", the handback timestamp is refreshed,
# and the stale "not the latest" error detail is cleared.

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("K3").Value = "2016-08-26 20:47:57"
$wsZhCn.Range("P3").Value = ""

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("K3").Value = "2016-08-26 20:48:10"
$wsDeDe.Range("P3").Value = ""

$wsZhCn.Columns.Item(16).ColumnWidth = 13.7470528738839
$wsDeDe.Columns.Item(16).ColumnWidth = 13.7470528738839
